$d = $word.ActiveDocument

$replacements = @(
    @("517÷7=", "605÷4="),
    @("844÷7=", "178÷3="),
    @("860÷5=", "759÷7="),
    @("615÷5=", "751÷9="),
    @("927÷8=", "198÷2="),
    @("649÷2=", "234÷9="),
    @("153÷2=", "220÷6="),
    @("331÷7=", "355÷4="),
    @("446÷2=", "842÷2="),
    @("464÷8=", "514÷2="),
    @("951÷3=", "728÷9="),
    @("745÷8=", "248÷8="),
    @("207÷2=", "986÷2="),
    @("445÷5=", "250÷2="),
    @("126÷7=", "768÷8="),
    @("654÷8=", "115÷7="),
    @("762÷6=", "129÷2="),
    @("393÷2=", "321÷2="),
    @("394÷6=", "398÷2="),
    @("510÷7=", "569÷4="),
    @("688÷2=", "588÷2="),
    @("782÷4=", "758÷2="),
    @("585÷8=", "772÷6="),
    @("199÷9=", "969÷3="),
    @("407÷7=", "174÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
